# Aggiornamento File Test Indicatori ISPRO
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ANALYSIS_UNIT")
$ws.Activate()

# --- Row 114 (style 28, unchanged) ---
$ws.Range("A114").Value = "ok"
$ws.Range("B114").Value = 2.5862069999999999
$ws.Range("D114").Value = "EST000003"
$ws.Range("G114").Value = 9

# --- Rows 115 & 116 need re-styling to style 25 (same as e.g. row 60) ---
$ws.Range("A60:H60").Copy()
$ws.Range("A115:H115").PasteSpecial(-4122)
$ws.Range("A60:H60").Copy()
$ws.Range("A116:H116").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A115").Value = "ko"
$ws.Range("D115").Value = "EST000019"
$ws.Range("G115").Value = 0

$ws.Range("A116").Value = "ko"
$ws.Range("D116").Value = "EST000007"
$ws.Range("G116").Value = 0

# --- Row 117 (style 1 / 27, unchanged) ---
$ws.Range("A117").Value = "ok"
$ws.Range("D117").Value = "EST000008"
$ws.Range("G117").Value = 0.2368421
$ws.Range("H117").Value = "E0001"

# --- Row 118 (style 1 / 27, unchanged) ---
$ws.Range("A118").Value = "ok"
$ws.Range("D118").Value = "EST000041"
$ws.Range("G118").Value = 0.2368421
$ws.Range("H118").Value = "E0001"

# --- Row 119 (style 1 / 27, unchanged) ---
$ws.Range("A119").Value = "ok"
$ws.Range("D119").Value = "EST000029"
$ws.Range("G119").Value = 0.2368421
$ws.Range("H119").Value = "E0001"

# --- Rows 120 & 121 need re-styling to style 25 (same as e.g. row 61/62) ---
$ws.Range("A61:H61").Copy()
$ws.Range("A120:H120").PasteSpecial(-4122)
$ws.Range("A61:H61").Copy()
$ws.Range("A121:H121").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A120").Value = "ko"
$ws.Range("D120").Value = "EST000013"
$ws.Range("G120").Value = 0

$ws.Range("A121").Value = "ko"
$ws.Range("D121").Value = "EST000018"
$ws.Range("G121").Value = 0

# --- Row 122 (style 1 / 27, unchanged) ---
$ws.Range("A122").Value = "ok"
$ws.Range("D122").Value = "EST000021"
$ws.Range("G122").Value = 0.2368421
$ws.Range("H122").Value = "E0001"

# --- Row 123 (style 1 / 27, unchanged) ---
$ws.Range("A123").Value = "ok"
$ws.Range("D123").Value = "EST000034"
$ws.Range("G123").Value = 0.2368421
$ws.Range("H123").Value = "E0001"

# --- Row 124 (style 1 / 27, unchanged) ---
$ws.Range("A124").Value = "ok"
$ws.Range("D124").Value = "EST000010"
$ws.Range("G124").Value = 0.2368421
$ws.Range("H124").Value = "E0001"

# --- Update sheet view: scroll position & selection ---
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("I127").Select()
